$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block spans rows 2-10 and columns A-I (9 columns).
# Turn it into an identity matrix: for row r (2..10), set a 1 in the
# column whose index equals (r-1), and 0 everywhere else in that row.
for ($r = 2; $r -le 10; $r++) {
    $diagCol = $r - 1
    for ($c = 1; $c -le 9; $c++) {
        if ($c -eq $diagCol) {
            $ws.Cells.Item($r, $c).Value = 1
        } else {
            $ws.Cells.Item($r, $c).Value = 0
        }
    }
}
